$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 was a blank templated row; copy the formatting from row 23 (the
# last filled-in row) down onto row 24 so the new entry matches the
# existing journal styling (date / time / duration number formats).
$ws.Range("A23:D23").Copy()
$ws.Range("A24:D24").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new journal entry.
$ws.Range("A24").Value = 43923
$ws.Range("B24").Value = 0.5625
$ws.Range("C24").Value = 0.62847222222222221
$ws.Range("D24").Formula = "=C24-B24"
$ws.Range("E24").Value = "Début de la rédaction de la documentation du projet"

# Match the author's final selection/cursor position.
[void]$ws.Range("E24").Select()
